$wb = $excel.ActiveWorkbook

# Fix typos in the shared member-title strings.
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("C2:C16").Value = "Coordinator"

$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("C3").Value = "Technical Secretary"

# Move the active/selected tab from Sheet3 to Sheet4, updating each
# sheet's remembered selection along the way.
$ws3.Activate()
$ws3.Range("D20").Select()

$ws4.Activate()
$ws4.Range("D24").Select()
